# Depersonalization script
# Replaces the original (identifying) patient id values in column A with a
# simple sequential id, and refreshes the accompanying diagnosis codes in
# column B with their depersonalized replacements.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New, depersonalized diagnosis codes (column B) for rows 2..69, in order.
$newDiagnosis = @(
    1,2,2,1,1,1,1,2,2,1,1,1,1,1,1,2,2,1,1,1,
    1,1,1,1,1,2,2,1,1,2,2,1,1,1,2,2,2,2,1,2,
    2,2,2,1,1,1,1,2,2,1,1,1,1,1,1,1,2,2,2,2,
    1,2,2,2,2,1,1,2
)

$firstRow = 2
$lastRow  = 69
$newId    = 137

for ($row = $firstRow; $row -le $lastRow; $row++) {
    # Sequential, anonymized identifier replacing the original value.
    $ws.Range("A$row").Value = $newId

    # Depersonalized diagnosis code.
    $ws.Range("B$row").Value = $newDiagnosis[$row - $firstRow]

    $newId++
}
